$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that previously sat in the empty
#    paragraph right before the "Design Decisions" heading.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Re-create "_GoBack" as a collapsed bookmark right after the sentence
#    "The Checkout page was modified to use the dummy payment system."
#    A zero-length Range cannot be handed directly to Bookmarks.Add, so a
#    throw-away placeholder character is inserted, the bookmark is wrapped
#    around it, and the placeholder text is then cleared - leaving a
#    collapsed bookmark exactly where the real edit placed it.
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("The Checkout page was modified to use the dummy payment system.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $find.Collapse(0)
    $find.InsertAfter([char]1) | Out-Null
    $find.MoveEnd(1, 1) | Out-Null
    $d.Bookmarks.Add("_GoBack", $find)
    $gb = $d.Bookmarks.Item("_GoBack")
    $gbr = $gb.Range
    $gbr.Text = ""
}

# ---------------------------------------------------------------------------
# 3) Append two new bibliography entries (same list style/numbering as the
#    surrounding "References" bullets) documenting code sourced online.
# ---------------------------------------------------------------------------
$newUrls = @(
    "https://stackoverflow.com/questions/730268/unique-random-string-generation",
    "https://docs.microsoft.com/en-us/dotnet/api/system.security.cryptography.md5?view=netframework-4.8"
)

# Locate the last "References" bullet (the Books_Flat_Icon_Vector.svg entry)
# by its known index, found once via Find, then walk forward paragraph by
# paragraph so each new bullet inherits the ListParagraph/numId=2 formatting.
$anchorParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Books_Flat_Icon_Vector.svg*") {
        $anchorParaIndex = $i
    }
}

foreach ($url in $newUrls) {
    $anchorRange = $d.Paragraphs($anchorParaIndex).Range
    $anchorRange.Collapse(0)
    $anchorRange.InsertParagraphAfter()
    $anchorParaIndex = $anchorParaIndex + 1
    $newRange = $d.Paragraphs($anchorParaIndex).Range
    $newRange.Text = $url
}
